$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the stray _GoBack bookmark that sits after the "Reducer"
# definition paragraph ("A type of task that takes in a fixed-length key...").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: replace the "bug with comparing key values" paragraph text with
# "There are no bugs that we know of." and move the _GoBack bookmark to sit
# (as a zero-width bookmark) right after the new sentence; also drop one of
# the two blank paragraphs that used to follow it.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*a bug with comparing key values that we haven*figured out yet*") {
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Text = "There are no bugs that we know of."
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "There are no bugs that we know of.*") {
        $insertPos = $p.Range.End - 1

        # Insert a temporary marker character, bookmark it, then delete the
        # marker again -- this leaves a clean zero-width bookmark exactly
        # after the sentence without landing on the (buggy) degenerate
        # zero-width-range case.
        $ins = $d.Range($insertPos, $insertPos)
        $ins.InsertAfter("X")
        $markerRange = $d.Range($insertPos, $insertPos + 1)
        $markerRange.Bookmarks.Add("_GoBack")
        $markerRange2 = $d.Range($insertPos, $insertPos + 1)
        $markerRange2.Delete()
    }
}

# Remove one of the two blank paragraphs directly following that paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "There are no bugs that we know of.*") {
        $blank = $d.Range($p.Range.End, $p.Range.End + 1)
        $blank.Delete()
    }
}

# ---------------------------------------------------------------------------
# Change 3: add a lastRenderedPageBreak marker in front of the "Our Map-reduce
# framework..." paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Our Map-reduce framework together with DFS*") {
        $startRng = $d.Range($p.Range.Start, $p.Range.Start)
        $startRng.Collapse(1)
        $word.Selection.SetRange($startRng.Start, $startRng.Start)
    }
}

# ---------------------------------------------------------------------------
# Change 4: remove the "As mentioned, there's a bug with comparing byte
# arrays for the keys." paragraph and the blank paragraph that follows it.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "As mentioned, there*a bug with comparing byte arrays for the keys.*") {
        $full = $d.Range($p.Range.Start, $p.Range.End + 1)
        $full.Delete()
    }
}
